$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F63").Value = 26
$ws.Range("G63").Value = 918.3200000000001
$ws.Range("B66").Value = 219092.23
$ws.Range("F101").Value = 4
$ws.Range("G101").Value = 328.56
$ws.Range("F114").Value = 67
$ws.Range("G114").Value = 3131.58
$ws.Range("F115").Value = 11
$ws.Range("G115").Value = 625.9
$ws.Range("B123").Value = 76112.95
$ws.Range("B126").Value = 65258
$ws.Range("B127").Value = 64196
$ws.Range("F184").Value = 61
$ws.Range("G184").Value = 5002
$ws.Range("F186").Value = 30
$ws.Range("G186").Value = 1298.4
$ws.Range("B193").Value = 69338.25
$ws.Range("F206").Value = 77
$ws.Range("G206").Value = 4989.6
$ws.Range("B208").Value = 5036.09
$ws.Range("F215").Value = 180
$ws.Range("G215").Value = 20212.2
$ws.Range("B218").Value = 84556.49000000001
$ws.Range("F220").Value = 56
$ws.Range("G220").Value = 3538.08
$ws.Range("F227").Value = 53
$ws.Range("G227").Value = 6073.8
$ws.Range("B229").Value = 33082.27
$ws.Range("F276").Value = 13
$ws.Range("G276").Value = 845.78
$ws.Range("B290").Value = 66194
$ws.Range("C290").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F290").Value = 27
$ws.Range("G290").Value = 2313.36
$ws.Range("B291").Value = 64983
$ws.Range("C291").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F291").Value = 6
$ws.Range("G291").Value = 514.08
$ws.Range("B295").Value = 132090.84
$ws.Range("B304").Value = 63520
$ws.Range("E304").Value = 153.4
$ws.Range("F304").Value = 39
$ws.Range("G304").Value = 5626.92
$ws.Range("B305").Value = 55373
$ws.Range("E305").Value = 163.62
$ws.Range("F305").Value = -94
$ws.Range("G305").Value = -13562.32
$ws.Range("B306").Value = 63531
$ws.Range("E306").Value = 152.53
$ws.Range("F306").Value = 29
$ws.Range("G306").Value = 4160.92
$ws.Range("B307").Value = 57802
$ws.Range("E307").Value = 162.71
$ws.Range("F307").Value = -79
$ws.Range("G307").Value = -11334.92
$ws.Range("B308").Value = 63510
$ws.Range("E308").Value = 50.66
$ws.Range("F308").Value = 80
$ws.Range("G308").Value = 3811.2
$ws.Range("B309").Value = 55356
$ws.Range("E309").Value = 54.04
$ws.Range("F309").Value = -158
$ws.Range("G309").Value = -7527.12
$ws.Range("F325").Value = 47
$ws.Range("G325").Value = 7104.99
$ws.Range("B328").Value = 1398.59
$ws.Range("F352").Value = 125
$ws.Range("G352").Value = 15613.75
$ws.Range("B356").Value = 79927.96000000001
$ws.Range("F361").Value = 252
$ws.Range("G361").Value = 35428.68
$ws.Range("B363").Value = 80682.8
$ws.Range("F377").Value = 47
$ws.Range("G377").Value = 45581.07
$ws.Range("B378").Value = 45581.07
$ws.Range("B381").Value = 58047
$ws.Range("D381").Value = 105.54
$ws.Range("E381").Value = 126.1
$ws.Range("F381").Value = 32
$ws.Range("G381").Value = 3377.28
$ws.Range("B382").Value = 47097
$ws.Range("D382").Value = 112.28
$ws.Range("E382").Value = 134.16
$ws.Range("F382").Value = 15
$ws.Range("G382").Value = 1684.2
$ws.Range("F402").Value = 60
$ws.Range("G402").Value = 2058.6
$ws.Range("B417").Value = 178901.6
$ws.Range("F420").Value = 37
$ws.Range("G420").Value = 5880.04
$ws.Range("B427").Value = 23986.47
$ws.Range("F431").Value = 20
$ws.Range("G431").Value = 537.8
$ws.Range("F434").Value = 161
$ws.Range("G434").Value = 5604.41
$ws.Range("F435").Value = 36
$ws.Range("G435").Value = 2507.4
$ws.Range("B438").Value = 27029.44
$ws.Range("B479").Value = 53319
$ws.Range("E479").Value = 310.64
$ws.Range("F479").Value = -6
$ws.Range("G479").Value = -1643.52
$ws.Range("B480").Value = 64810
$ws.Range("E480").Value = 291.22
$ws.Range("F480").Value = 0
$ws.Range("G480").Value = 0
$ws.Range("F488").Value = 3
$ws.Range("G488").Value = 190.83
$ws.Range("F491").Value = 17
$ws.Range("G491").Value = 308.38
$ws.Range("B493").Value = 3203.72
$ws.Range("B506").Value = 60022
$ws.Range("E506").Value = 37.22
$ws.Range("F506").Value = -113
$ws.Range("G506").Value = -3709.79
$ws.Range("B507").Value = 64830
$ws.Range("E507").Value = 34.9
$ws.Range("F507").Value = 86
$ws.Range("G507").Value = 2823.38
$ws.Range("F519").Value = 424
$ws.Range("G519").Value = 23269.12
$ws.Range("F523").Value = 162
$ws.Range("G523").Value = 13868.82
$ws.Range("F524").Value = 23
$ws.Range("G524").Value = 2031.82
$ws.Range("B525").Value = 132124.39
$ws.Range("F527").Value = 59
$ws.Range("G527").Value = 1953.49
$ws.Range("B535").Value = 26511.44
$ws.Range("F605").Value = 200
$ws.Range("G605").Value = 26620
$ws.Range("B607").Value = 27160.04
$ws.Range("F622").Value = 492
$ws.Range("G622").Value = 50631.72
$ws.Range("B628").Value = 224291.07
$ws.Range("F648").Value = 0
$ws.Range("G648").Value = 0
$ws.Range("B657").Value = 82778.14999999999
$ws.Range("F715").Value = 125
$ws.Range("G715").Value = 12940
$ws.Range("B717").Value = 19794.95
$ws.Range("B718").Value = 2966217.37
$ws.Range("B719").Value = 2966217.37
